$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("G4").Value = 1.57
$ws.Range("H4").Value = 3.6
$ws.Range("J4").Value = 2.2
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("AA4").Value = 15
$ws.Range("AB4").Value = 34
$ws.Range("AC4").Value = 8.5
$ws.Range("AF4").Value = 67
$ws.Range("AH4").Value = 13
$ws.Range("AJ4").Value = 19
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 26
$ws.Range("AS4").Value = 201
$ws.Range("AT4").Value = 2.63
$ws.Range("AU4").Value = 9.5
$ws.Range("AV4").Value = 67
$ws.Range("AZ4").Value = 126
$ws.Range("BA4").Value = 151

# Row 5
$ws.Range("G5").Value = 4.5
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 1.95
$ws.Range("K5").Value = 1.91
$ws.Range("L5").Value = 2.75
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.38
$ws.Range("Q5").Value = 2.7
$ws.Range("R5").Value = 1.44
$ws.Range("S5").Value = 1.57
$ws.Range("T5").Value = 2.25
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("W5").Value = 9
$ws.Range("X5").Value = 21
$ws.Range("AC5").Value = 6
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 81
$ws.Range("AH5").Value = 5.5
$ws.Range("AJ5").Value = 9.5
$ws.Range("AK5").Value = 17
$ws.Range("AL5").Value = 21
$ws.Range("AT5").Value = 2.25
$ws.Range("AU5").Value = 9.5
$ws.Range("AX5").Value = 12
$ws.Range("AY5").Value = 29
$ws.Range("BB5").Value = 251

# Row 6
$ws.Range("AG6").Value = 700

# Row 12
$ws.Range("G12").Value = 1.23
$ws.Range("H12").Value = 6.2
$ws.Range("I12").Value = 9.75
$ws.Range("J12").Value = 1.57
$ws.Range("K12").Value = 2.95
$ws.Range("L12").Value = 7.4
$ws.Range("S12").Value = 1.2
$ws.Range("T12").Value = 4
$ws.Range("U12").Value = 1.78
$ws.Range("V12").Value = 1.93
$ws.Range("W12").Value = 11
$ws.Range("X12").Value = 7.9
$ws.Range("Z12").Value = 8.25
$ws.Range("AC12").Value = 10.5
$ws.Range("AH12").Value = 35
$ws.Range("AI12").Value = 80
$ws.Range("AJ12").Value = 30
$ws.Range("AK12").Value = 250
$ws.Range("AL12").Value = 100
$ws.Range("AM12").Value = 75
$ws.Range("AN12").Value = 3.45
$ws.Range("AP12").Value = 13
$ws.Range("AT12").Value = 4
$ws.Range("AU12").Value = 8.5
$ws.Range("AV12").Value = 60
$ws.Range("AW12").Value = 10.5
$ws.Range("AX12").Value = 50
$ws.Range("AZ12").Value = 350
$ws.Range("BA12").Value = 250
